$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("409:409").Insert()

$ws.Range("A409").Value = 11
$ws.Range("B409").Value = "Vega Monumental Concepción"
$ws.Range("C409").Value = "Bíobío"
$ws.Range("D409").Value = 45209
$ws.Range("E409").Value = 8
$ws.Range("F409").Value = 100114013
$ws.Range("G409").Value = "Zanahoria"
$ws.Range("H409").Value = "Sin especificar"
$ws.Range("I409").Value = "Primera"
$ws.Range("J409").Value = 160
$ws.Range("K409").Value = 5000
$ws.Range("L409").Value = 5500
$ws.Range("M409").Value = 5250
$ws.Range("N409").Value = "$/saco 20 kilos"
$ws.Range("O409").Value = "Región de Ñuble"
$ws.Range("P409").Value = 262
$ws.Range("Q409").Value = 20
$ws.Range("R409").Value = "Hortaliza"
